$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.527.47'
$ws.Range('E2').Value = '  -0.28%  '

$ws.Range('D3').Value = '1.807.12'
$ws.Range('E3').Value = '  -0.51%  '

$ws.Range('E5').Value = '  -0.02%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.609'
$ws.Range('E6').Value = '  +8.69%  '

$ws.Range('E7').Value = '  +0.13%  '

$ws.Range('B8').Value = 'WrappedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D8').Value = '5.696.75'
$ws.Range('E8').Value = '  +213.58%  '

$ws.Range('B9').Value = 'Solana'
$ws.Range('C9').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '36.60'
$ws.Range('E9').Value = '  +5.11%  '

$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.302'
$ws.Range('E10').Value = '  +0.76%  '

$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0696'
$ws.Range('E11').Value = '  -0.02%  '

$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0965'
$ws.Range('E12').Value = '  +1.33%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.067.45'
$ws.Range('E13').Value = '  -0.55%  '

$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.51'
$ws.Range('E14').Value = '  +1.61%  '

$ws.Range('E15').Value = '  +1.62%  '

$ws.Range('E16').Value = '  +3.78%  '

$ws.Range('D17').Value = '34.484.86'
$ws.Range('E17').Value = '  -0.48%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.94'
$ws.Range('E18').Value = '  +1.10%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.22'
$ws.Range('E19').Value = '  -0.59%  '

$ws.Range('E20').Value = '  -1.31%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.60'
$ws.Range('E21').Value = '  +0.31%  '

$ws.Range('E22').Value = '  +0.12%  '

$ws.Range('E23').Value = '  -0.64%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.22'
$ws.Range('E24').Value = '  +6.19%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.40'
$ws.Range('E25').Value = '  -0.42%  '

$ws.Range('E26').Value = '  +6.64%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.43'
$ws.Range('E27').Value = '  +3.89%  '

$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.123'
$ws.Range('E28').Value = '  +4.93%  '

$ws.Range('E30').Value = '  -1.51%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0531'
$ws.Range('E31').Value = '  -0.43%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.85'
$ws.Range('E32').Value = '  -0.21%  '

$ws.Range('E33').Value = '  -0.19%  '

$ws.Range('E34').Value = '  -1.87%  '

$ws.Range('D35').Value = '1.394.28'
$ws.Range('E35').Value = '  -1.71%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.672'
$ws.Range('E36').Value = '  -0.57%  '

$ws.Range('E37').Value = '  -5.68%  '

$ws.Range('E38').Value = '  -0.25%  '

$ws.Range('E39').Value = '  -1.13%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.967'
$ws.Range('E40').Value = '  +0.64%  '

$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '82.57'
$ws.Range('E41').Value = '  -4.04%  '

$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.83'
$ws.Range('E42').Value = '  -0.83%  '

$ws.Range('E43').Value = '  +0.83%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.20'
$ws.Range('E44').Value = '  +7.97%  '

$ws.Range('E45').Value = '  -2.75%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.03'
$ws.Range('E46').Value = '  -1.25%  '

$ws.Range('B47').Value = 'Kaspa'
$ws.Range('C47').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0495'
$ws.Range('E47').Value = '  -5.92%  '

$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '1.968.53'
$ws.Range('E48').Value = '  -0.59%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.28'
$ws.Range('E49').Value = '  -1.70%  '

$ws.Range('E50').Value = '  +0.15%  '

$ws.Range('D51').Value = '0.0₆0128'
$ws.Range('E51').Value = '  -2.97%  '
